$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---------------------------------------------------
# C (Edad) and D (Modalidad) get narrower, new E (Peso) column is added.
# NOTE: this engine's stored OOXML <col width> ends up as
# (COM ColumnWidth + 5/6), matching Excel's own char-width padding
# quirk, so back the input off by 5/6 to land exactly on 6 / 11 / 7.
$ws.Columns.Item(3).ColumnWidth = 6 - 5/6
$ws.Columns.Item(4).ColumnWidth = 11 - 5/6
$ws.Columns.Item(5).ColumnWidth = 7 - 5/6

# --- New column E: copy the header/data formatting from column C -----
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

# --- Header row --------------------------------------------------------
$ws.Range("C1").Value = "Edad"
$ws.Range("D1").Value = "Modalidad"
$ws.Range("E1").Value = "Peso"

# --- Row 2 (Paola) ------------------------------------------------------
$ws.Range("C2").Value = "'25"
$ws.Range("D2").Value = "Boxeo"
$ws.Range("E2").Value = "54 kg"

# --- Row 3 (now AngeloO Lopez) ------------------------------------------
$ws.Range("A3").Value = "'1750321901"
$ws.Range("B3").Value = "AngeloO Lopez"
$ws.Range("C3").Value = "'21"
$ws.Range("D3").Value = "Boxeo"
$ws.Range("E3").Value = "58 kg"
